$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: fill in the Score (F7) that was missing - this game's score is "B+"
$ws.Range("F7").Value2 = "B+"

# Row 8 was a duplicate entry of row 7's game; correct it with what used to be
# row 9's data (shifted up), fixing the "same game added twice" mistake.
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("D8").Value2 = $ws.Range("D9").Value2
$ws.Range("F8").Value2 = "A+"
$ws.Range("J8").Value2 = 208

# Row 9 becomes a new entry (replacing the old duplicate-derived one)
$ws.Range("D9").NumberFormat = "h:mm"
$ws.Range("D9").Value2 = 0.97986111111111107
$ws.Range("F9").Value2 = "A-"
$ws.Range("J9").Value2 = 162

# Update the selected cell to E11 (per the sheetView selection change)
$ws.Range("E11").Select() | Out-Null
